# "Generate Report for Handback"
#
# This localization-status workbook tracks per-file handoff/handback state
# for each target locale (zh-cn, de-de). Running the handback report:
#   1. Marks the files as handed back / in sync with en-US (Status column).
#   2. Fills in the "Latest Target File", "Latest Handback File" and
#      "Latest Handback DateTime" columns on the per-locale sheets, turning
#      the target-file cell into a hyperlink back to the source doc (just
#      like the existing "Source File Name" column does).
#   3. Widens a few columns that now hold longer file names / datetimes.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

$mdBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/06b8b5fed0864774f689490a00885d9a7d5f693e/e2e/"

$file1Name = "2db4db43-ed49-4db6-94a7-647e8cb93e42.md"
$file2Name = "dccdc1f8-1fc0-4ee0-ac3a-7fe9156b3bc5.md"

# ---------------------------------------------------------------------
# Overview sheet: update the per-locale status text shown for each file.
# ---------------------------------------------------------------------
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# Widen the zh-cn / de-de status columns to fit the longer text.
$overview.Columns.Item(5).ColumnWidth = 29.14
$overview.Columns.Item(6).ColumnWidth = 29.14

# ---------------------------------------------------------------------
# zh-cn sheet (row2 = file1, row3 = file2)
# ---------------------------------------------------------------------
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), ($mdBase + $file1Name), "", "", $file1Name)
$zhcn.Range("J2").Value = "2db4db43-ed49-4db6-94a7-647e8cb93e42.7c2217b3072ade006b225d85acc5ed30726c04e1.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-30 11:10:55"

$zhcn.Hyperlinks.Add($zhcn.Range("I3"), ($mdBase + $file2Name), "", "", $file2Name)
$zhcn.Range("J3").Value = "dccdc1f8-1fc0-4ee0-ac3a-7fe9156b3bc5.5e35570a73f1e5e1d75101395355b124d563c1cf.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-30 11:10:55"

# Match the hyperlink look already used by column A.
$zhcn.Range("I2").Font.Underline = 2
$zhcn.Range("I2").Font.Color = 15570276
$zhcn.Range("I3").Font.Underline = 2
$zhcn.Range("I3").Font.Color = 15570276

$zhcn.Columns.Item(3).ColumnWidth = 29.14
$zhcn.Columns.Item(9).ColumnWidth = 39.14
$zhcn.Columns.Item(10).ColumnWidth = 39.14

# ---------------------------------------------------------------------
# de-de sheet (row2 = file1, row3 = file2)
# ---------------------------------------------------------------------
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

$dede.Hyperlinks.Add($dede.Range("I2"), ($mdBase + $file1Name), "", "", $file1Name)
$dede.Range("J2").Value = "2db4db43-ed49-4db6-94a7-647e8cb93e42.7c2217b3072ade006b225d85acc5ed30726c04e1.de-de.xlf"
$dede.Range("K2").Value = "2016-08-30 11:11:07"

$dede.Hyperlinks.Add($dede.Range("I3"), ($mdBase + $file2Name), "", "", $file2Name)
$dede.Range("J3").Value = "dccdc1f8-1fc0-4ee0-ac3a-7fe9156b3bc5.5e35570a73f1e5e1d75101395355b124d563c1cf.de-de.xlf"
$dede.Range("K3").Value = "2016-08-30 11:11:07"

$dede.Range("I2").Font.Underline = 2
$dede.Range("I2").Font.Color = 15570276
$dede.Range("I3").Font.Underline = 2
$dede.Range("I3").Font.Color = 15570276

$dede.Columns.Item(3).ColumnWidth = 29.14
$dede.Columns.Item(9).ColumnWidth = 39.14
$dede.Columns.Item(10).ColumnWidth = 39.14

Write-Host "Handback report generated."
